$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Vishay Dale 100 ohm SMD resistor -> Stackpole Electronics 330 ohm through-hole resistor
$ws.Range("A7").Value = "Stackpole Electronics RSMF2JT330R "
$ws.Range("B7").Value = "Through Hole 330 ohm 2 watt resistor"
$ws.Range("C7").Value = 0.29
$ws.Range("G7").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/RSF2JT330R/RSF2JT330RCT-ND/2021796"

# Row 8: Vishay Dale 10 ohm SMD resistor -> Stackpole Electronics 10 ohm through-hole resistor
$ws.Range("A8").Value = "Stackpole Electronics  CF12JT10R0 "
$ws.Range("B8").Value = "Through Hole 10 ohm ½ watt Resistor"
$ws.Range("C8").Value = 0.1
$ws.Range("G8").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF12JT10R0/CF12JT10R0CT-ND/1830446"

# Update the active cell / selection to reflect the author's cursor position after editing
$ws.Range("B18").Select()
